# Generate Report for Handoff
#
# The localization-status report was regenerated and the handoff file
# "62d7af38-5806-4a20-864e-397ece57a52c.md" picked up a fresh "Ready for
# handoff" timestamp. Update the three places that record that timestamp:
#   - Overview!G6            "Latest HO Xliff Generate Date"
#   - zh-cn!H6                "Latest Handoff Datetime"
#   - de-de!H6                "Latest Handoff Datetime"

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G6").Value = "2016-08-27 02:41:05"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H6").Value = "2016-08-27 02:40:56"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H6").Value = "2016-08-27 02:41:05"
